# Generate Report for Handoff
# Adds two more "rows" of handoff/report data (for the two new source files
# b25b6e5f-...md and ff761548-...png) to every sheet, alongside refreshed
# values for the already-tracked 46e935ec-...png (previously e5102ddf-...md)
# entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2: refresh the existing tracked file's name + latest handoff date.
$ov.Range("A2").Hyperlinks.Delete()
$ov.Range("A2").Value = "46e935ec-b109-44d5-b61d-49b37aab4348.png"
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a64f29898e0884a58eca8e421dfa80bf1f0ce03c/e2e/46e935ec-b109-44d5-b61d-49b37aab4348.png", "", "", "46e935ec-b109-44d5-b61d-49b37aab4348.png")
$ov.Range("A2").Style = "HyperLink"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-52-18 16:52:50"

# Row 3: new file b25b6e5f-...md
$ov.Range("A3").Value = "b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.md"
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/a64f29898e0884a58eca8e421dfa80bf1f0ce03c/e2e/b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.md", "", "", "b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.md")
$ov.Range("A3").Style = "HyperLink"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-52-18 16:52:50"

# Row 4: new file ff761548-...png
$ov.Range("A4").Value = "ff761548-4dcd-45ae-aaf0-3893c2453015.png"
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/a64f29898e0884a58eca8e421dfa80bf1f0ce03c/e2e/ff761548-4dcd-45ae-aaf0-3893c2453015.png", "", "", "ff761548-4dcd-45ae-aaf0-3893c2453015.png")
$ov.Range("A4").Style = "HyperLink"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-52-18 16:52:50"

# ---------------------------------------------------------------------------
# Per-language detail sheets ("zh-cn" and "de-de")
# ---------------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; DateTime = "2016-03-18 16:52:47"; RowDateTime = "2016-03-18 16:52:47";
       HandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea6876ea89b93bd47d99707ea2f32f1ab9eebe6b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/";
       XlfFile = "b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.fc4cd1a358e4b7b50dd0f520ec612def86e9b149.zh-cn.xlf" },
    @{ Sheet = "de-de"; DateTime = "2016-03-18 16:52:50"; RowDateTime = "2016-03-18 16:52:50";
       HandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d043f5ae3aeb0dc5d57d7eef4f6888a2988ba9ec/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/";
       XlfFile = "b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.fc4cd1a358e4b7b50dd0f520ec612def86e9b149.de-de.xlf" }
)

$srcBase = "https://github.com/OpenLocalizationTest/oltest/blob/a64f29898e0884a58eca8e421dfa80bf1f0ce03c/e2e/"

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # --- Row 2: refresh existing "46e935ec-...png" (was e5102ddf-...md) ---
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("A2").Value = "46e935ec-b109-44d5-b61d-49b37aab4348.png"
    $ws.Hyperlinks.Add($ws.Range("A2"), ($srcBase + "46e935ec-b109-44d5-b61d-49b37aab4348.png"), "", "", "46e935ec-b109-44d5-b61d-49b37aab4348.png")
    $ws.Range("A2").Style = "HyperLink"

    $ws.Range("B2").Hyperlinks.Delete()
    $ws.Range("B2").Value = ".png"
    $ws.Hyperlinks.Add($ws.Range("B2"), ($srcBase + "46e935ec-b109-44d5-b61d-49b37aab4348.png"), "", "", ".png")
    $ws.Range("B2").Style = "HyperLink"

    $ws.Range("C2").Value = "Ready for handoff"

    $ws.Range("D2").Hyperlinks.Delete()
    $ws.Range("D2").Value = "763d60acbfd46075c1e473634cca645e6cdab3ed.png"
    $ws.Hyperlinks.Add($ws.Range("D2"), ($lang.HandoffBase + "763d60acbfd46075c1e473634cca645e6cdab3ed.png"), "", "", "763d60acbfd46075c1e473634cca645e6cdab3ed.png")
    $ws.Range("D2").Style = "HyperLink"

    $ws.Range("E2").Value = $lang.DateTime
    $ws.Range("H2").Value = "0001-01-01 00:00:00"
    $ws.Range("I2").Value = "IsDependency"
    $ws.Range("J2").Value = 'e2e\b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.md'

    # --- Row 3: new file b25b6e5f-...md ---
    $ws.Range("A3").Value = "b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.md"
    $ws.Hyperlinks.Add($ws.Range("A3"), ($srcBase + "b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.md"), "", "", "b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.md")
    $ws.Range("A3").Style = "HyperLink"

    $ws.Range("B3").Value = ".md"
    $ws.Hyperlinks.Add($ws.Range("B3"), ($srcBase + "b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.md"), "", "", ".md")
    $ws.Range("B3").Style = "HyperLink"

    $ws.Range("C3").Value = "Ready for handoff"

    $ws.Range("D3").Value = $lang.XlfFile
    $ws.Hyperlinks.Add($ws.Range("D3"), ($lang.HandoffBase + $lang.XlfFile), "", "", $lang.XlfFile)
    $ws.Range("D3").Style = "HyperLink"

    $ws.Range("E3").Value = $lang.DateTime
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "Include"

    # --- Row 4: new file ff761548-...png ---
    $ws.Range("A4").Value = "ff761548-4dcd-45ae-aaf0-3893c2453015.png"
    $ws.Hyperlinks.Add($ws.Range("A4"), ($srcBase + "ff761548-4dcd-45ae-aaf0-3893c2453015.png"), "", "", "ff761548-4dcd-45ae-aaf0-3893c2453015.png")
    $ws.Range("A4").Style = "HyperLink"

    $ws.Range("B4").Value = ".png"
    $ws.Hyperlinks.Add($ws.Range("B4"), ($srcBase + "ff761548-4dcd-45ae-aaf0-3893c2453015.png"), "", "", ".png")
    $ws.Range("B4").Style = "HyperLink"

    $ws.Range("C4").Value = "Ready for handoff"

    $ws.Range("D4").Value = "001916a9f5d1536fa975967c830bdba72f801172.png"
    $ws.Hyperlinks.Add($ws.Range("D4"), ($lang.HandoffBase + "001916a9f5d1536fa975967c830bdba72f801172.png"), "", "", "001916a9f5d1536fa975967c830bdba72f801172.png")
    $ws.Range("D4").Style = "HyperLink"

    $ws.Range("E4").Value = $lang.DateTime
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "IsDependency"
    $ws.Range("J4").Value = 'e2e\b25b6e5f-cbc2-4511-bc1e-dd831d8076c9.md'
}

Write-Host "Handoff report rows added."
